$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "HoaiBao"
$ws.Range("B4").Value = "123456aA@"
